# Reorders the "readme" summary table's data columns from
#   index, JobNo, sheet_name, Author, Date
# to
#   index, Author, JobNo, sheet_name, Date
# (Author column moved to be the first data column), and bumps the
# "Date of Analysis" timestamp on the "Project Information" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")

# Capture the existing per-row text (B=JobNo, C=sheet_name, D=Author) before
# overwriting anything, since the new layout reuses these values. Use .Text
# (not .Value2) so numeric-looking strings like "/c/e" stay text and we don't
# accidentally coerce anything to a number. Column E (Date) is left untouched.
$lastRow = 12
$jobNoVals = @{}
$sheetNameVals = @{}
$authorVals = @{}

for ($r = 2; $r -le $lastRow; $r++) {
    $jobNoVals[$r] = $ws.Cells.Item($r, 2).Text
    $sheetNameVals[$r] = $ws.Cells.Item($r, 3).Text
    $authorVals[$r] = $ws.Cells.Item($r, 4).Text
}

# Re-header the table: B=Author, C=JobNo, D=sheet_name (E=Date unchanged)
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "JobNo"
$ws.Range("D1").Value = "sheet_name"

# Rewrite each data row in the new column order.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $authorVals[$r]
    $ws.Cells.Item($r, 3).Value = $jobNoVals[$r]
    $ws.Cells.Item($r, 4).Value = $sheetNameVals[$r]
}

# Bump the recorded analysis timestamp on "Project Information".
$wsInfo = $wb.Worksheets.Item("Project Information")
$wsInfo.Range("B12").Value = "2022-06-15 15:57:20.064427"
